$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# "dicionario de dados" update: add two new table blocks ("Post" and
# "like") below the existing "Parada" block, following the exact same
# visual layout/formatting pattern used by the other table blocks.
# ------------------------------------------------------------------

# --- Post block (rows 26-33) -------------------------------------
# Use the "Parada" block (rows 18-24) as the formatting template: it
# has the exact same style pattern (title row, column-header row, 5
# attribute rows) that the new "Post" block needs for its first 7
# rows.
$ws.Range("A18:C24").Copy($ws.Range("A26:C32"))
# Row 33 (extra FK attribute row) reuses the same row style as rows
# 21/23/29/31 (odd attribute rows).
$ws.Range("A23:C23").Copy($ws.Range("A33:C33"))

# --- like block (rows 35-38) --------------------------------------
$ws.Range("A18:C18").Copy($ws.Range("A35:C35"))
$ws.Range("A19:C19").Copy($ws.Range("A36:C36"))
$ws.Range("A24:C24").Copy($ws.Range("A37:C37"))
$ws.Range("A23:C23").Copy($ws.Range("A38:C38"))

# --- Row heights (match the other table blocks) --------------------
$ws.Range("A26:C26").RowHeight = 21
$ws.Range("A27:C27").RowHeight = 15.75
$ws.Range("A28:C33").RowHeight = 15.75
$ws.Range("A35:C35").RowHeight = 21
$ws.Range("A36:C36").RowHeight = 15.75
$ws.Range("A37:C38").RowHeight = 15.75

# --- Cell values -----------------------------------------------------
# NOTE: values are assigned in the same order the original author
# typed them in (column by column within each new block) so that the
# shared-strings table ends up with the same append order/indices.

# Post block - column A first (title, then each attribute name)
$ws.Range("A26").Value2 = "Post"
$ws.Range("A27").Value2 = "Nome do Atributo"
$ws.Range("A28").Value2 = "idPost"
$ws.Range("A29").Value2 = "titulo"
$ws.Range("A30").Value2 = "descricao"
$ws.Range("A31").Value2 = "likes"
$ws.Range("A32").Value2 = "Usuario_idUsuario"
$ws.Range("A33").Value2 = "Viagem_idViagem"

# Post block - column B (descriptions)
$ws.Range("B27").Value2 = "Dado a ser recebido"
$ws.Range("B28").Value2 = "Id fornecido para identificação do post"
$ws.Range("B29").Value2 = "Titulo do post postado pelo usuario"
$ws.Range("B30").Value2 = "Descrição do post fornecida pelo usuario"
$ws.Range("B31").Value2 = "Numero de likes que o post recebeu"
$ws.Range("B32").Value2 = "Chave estrangeira referente aos Usuarios e seus post's"
$ws.Range("B33").Value2 = "Chave estrangeira referente aos Post de viagens "

# Post block - column C (types)
$ws.Range("C27").Value2 = "Tipo de Dado"
$ws.Range("C28").Value2 = "Int"
$ws.Range("C29").Value2 = "Varchar"
$ws.Range("C30").Value2 = "Varchar"
$ws.Range("C31").Value2 = "Int"
$ws.Range("C32").Value2 = "Int"
$ws.Range("C33").Value2 = "Int"

# like block - column A
$ws.Range("A35").Value2 = "like"
$ws.Range("A36").Value2 = "Nome do Atributo"
$ws.Range("A37").Value2 = "Usuario_idUsuario"
$ws.Range("A38").Value2 = "Post_idPost"

# like block - column B
$ws.Range("B36").Value2 = "Dado a ser recebido"
$ws.Range("B37").Value2 = "Chave estrangeira referente aos Usuarios e seus post's"
$ws.Range("B38").Value2 = "Chave estrangeira referente aos posts e seu identificador"

# like block - column C
$ws.Range("C36").Value2 = "Tipo de Dado"
$ws.Range("C37").Value2 = "Int"
$ws.Range("C38").Value2 = "int"

# --- View state: scroll position / selection, matching the diff ----
$ws.Application.ActiveWindow.ScrollRow = 21
$ws.Range("D38").Select()

Write-Output "edit applied"
